$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new weekly records were added to the dataset. Inserting blank rows at
# 11 and 23 shifts the existing rows down (11->12, ..., 22->23->24 etc.),
# reproducing the row-shuffle seen in the diff without touching any of the
# values that merely moved to a new row.
$ws.Rows.Item(11).Insert()
$ws.Rows.Item(23).Insert()

# New row 11 (brand new data point, date 2021-07-06)
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C11").Value = "Metropolitana"
$ws.Range("D11").Value = 44383
$ws.Range("E11").Value = 13
$ws.Range("F11").Value = 100112035
$ws.Range("G11").Value = "Bruselas (repollito)"
$ws.Range("H11").Value = "Sin especificar"
$ws.Range("I11").Value = "Primera"
$ws.Range("J11").Value = 25
$ws.Range("K11").Value = 13000
$ws.Range("L11").Value = 14000
$ws.Range("M11").Value = 13480
$ws.Range("N11").Value = "`$/malla 15 kilos"
$ws.Range("O11").Value = "Hijuelas"
$ws.Range("P11").Value = 899
$ws.Range("Q11").Value = 15
$ws.Range("R11").Value = "Hortaliza"

# New row 23 (brand new data point, date 2022-06-07)
$ws.Range("A23").Value = 9
$ws.Range("B23").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C23").Value = "Metropolitana"
$ws.Range("D23").Value = 44719
$ws.Range("E23").Value = 13
$ws.Range("F23").Value = 100112035
$ws.Range("G23").Value = "Bruselas (repollito)"
$ws.Range("H23").Value = "Sin especificar"
$ws.Range("I23").Value = "Primera"
$ws.Range("J23").Value = 43
$ws.Range("K23").Value = 17000
$ws.Range("L23").Value = 18000
$ws.Range("M23").Value = 17512
$ws.Range("N23").Value = "`$/malla 15 kilos"
$ws.Range("O23").Value = "Hijuelas"
$ws.Range("P23").Value = 1167
$ws.Range("Q23").Value = 15
$ws.Range("R23").Value = "Hortaliza"

Write-Output "edit applied"
